# "Add files via upload" — update CUARTO SPRINT (4th sprint) backlog rows:
# mark the remaining tasks as completed, assign the people who did them,
# fill in actual hours, and fix the "FERNANDO FERNÁNDEZ" name accent
# wherever it appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the misspelled/unaccented name everywhere it is used so the shared
# string is corrected in place instead of leaving a duplicate behind.
$ws.Range("G13").Value = "FERNANDO FERNÁNDEZ"
$ws.Range("G28").Value = "FERNANDO FERNÁNDEZ"

# Sprint 4 summary row: now fully completed (160 of 160 hours).
$ws.Range("J41").Value = 1
$ws.Range("L41").Value = 160

# Row 42: "Implementar Funciones de Usuarios" -> Completado, by Nicólas Ramirez
$ws.Range("B42").Value = "Completado"
$ws.Range("G42").Value = "NICÓLAS RAMIREZ"
$ws.Range("J42").Value = 1
$ws.Range("L42").Value = 40

# Row 43: "Diseñar Vista de Usuarios" -> Completado, by Daniel Vicente
$ws.Range("B43").Value = "Completado"
$ws.Range("G43").Value = "DANIEL VICENTE"
$ws.Range("J43").Value = 1
$ws.Range("L43").Value = 40

# Row 44: "Implementar Controlador de Usuarios" -> Completado, by Álvaro López
$ws.Range("B44").Value = "Completado"
$ws.Range("G44").Value = "ÁLVARO LÓPEZ"
$ws.Range("J44").Value = 1
$ws.Range("L44").Value = 40

# Row 46: "Redactar READ.md" -> Completado, by Fernando Fernández
$ws.Range("B46").Value = "Completado"
$ws.Range("G46").Value = "FERNANDO FERNÁNDEZ"
$ws.Range("J46").Value = 1
$ws.Range("L46").Value = 40
